$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1544
$ws1.Range("F5").Value = 236
$ws1.Range("F7").Value = 836
$ws1.Range("F8").Value = 10053
$ws1.Range("F14").Value = 6966
$ws1.Range("F18").Value = 212

# Sheet "全部类型" (all types) — same events, different row positions
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1544
$ws4.Range("F5").Value = 236
$ws4.Range("F8").Value = 836
$ws4.Range("F11").Value = 10053
$ws4.Range("F17").Value = 6966
$ws4.Range("F21").Value = 212
